# Auto-upload VRF Excel file
# Append a new worksheet named "zxc" at the end of the workbook with the
# standard VRF header row (Outdoor/Indoor Model, Quantity, Serial(s)),
# matching the bold/bordered/centered header style used on the other
# sheets in this workbook.

$wb = $excel.ActiveWorkbook

# Remember the sheet that is active before we add the new one so the
# workbook's active-tab selection is left exactly as it was.
$prevActive = $wb.ActiveSheet

# Add the new sheet immediately after the last existing sheet, so it
# lands at the very end of the tab strip (sheetId 22 / last position).
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "zxc"

# Header row content.
$ws.Range("A1").Value = "Outdoor Model"
$ws.Range("B1").Value = "Outdoor Quantity"
$ws.Range("C1").Value = "Outdoor Serial(s)"
$ws.Range("D1").Value = "Indoor Model"
$ws.Range("E1").Value = "Indoor Quantity"
$ws.Range("F1").Value = "Indoor Serial(s)"

# Header formatting: bold, centered/top-aligned, thin box border - same
# look as the header row on every other sheet in this workbook.
$header = $ws.Range("A1:F1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108   # xlCenter
$header.VerticalAlignment = -4160     # xlTop
$header.Borders.LineStyle = 1         # xlContinuous (thin box border)

# Restore the originally active sheet/tab selection.
$prevActive.Activate()
